# Add a new "kill_equipped" mission-type tracker row to the
# missionTypeDefinitions table (Table13), right after the "kill_frozen"
# row (B62:H62), so it becomes the new row 63 and everything below it
# (including the other mission tables) shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at row 63 inside the missionTypeDefinitions table ---
$lo = $ws.ListObjects.Item("Table13")
$ws.Rows.Item(63).Insert()
$lo.Resize($ws.Range("B48:H67"))

# --- 2. Populate the new row's values ---
$ws.Range("B63").Value = "<Definition>"
$ws.Range("C63").Value = "kill_equipped"
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 7
$ws.Range("F63").Value = 1
$ws.Range("G63").Value = "TID_MISSION_OBJECTIVE_KILL_EQUIPABLE_DESC_SINGLE_RUN"
$ws.Range("H63").Value = "TID_MISSION_OBJECTIVE_KILL_EQUIPABLE_DESC_MULTI_RUN"

# --- 3. Match formatting of the row above, then touch up the borders so ---
#        the new (now "interior") row gets a full box border, matching
#        how Excel re-flows table borders when a row is inserted mid-table.
$ws.Range("B62:H62").Copy()
$ws.Range("B63:H63").PasteSpecial(-4122)

$rngB = $ws.Range("B63")
$rngB.Borders.Item(7).LineStyle = 0
$rngB.Borders.Item(10).LineStyle = 1
$rngB.Borders.Item(8).LineStyle = 1
$rngB.Borders.Item(9).LineStyle = 1

$rngCF = $ws.Range("C63:F63")
$rngCF.Borders.Item(7).LineStyle = 1
$rngCF.Borders.Item(10).LineStyle = 1
$rngCF.Borders.Item(8).LineStyle = 1
$rngCF.Borders.Item(9).LineStyle = 1

# --- 4. The tables below the insertion point don't auto-resize, fix them ---
$ws.ListObjects.Item("missionDifficultyDefinitions").Resize($ws.Range("B71:L74"))
$ws.ListObjects.Item("Table13303132").Resize($ws.Range("B79:E91"))
$ws.ListObjects.Item("Table1330313234").Resize($ws.Range("B95:E98"))
$ws.ListObjects.Item("Table133031323435").Resize($ws.Range("B102:D103"))

# --- 5. The duplicate-values conditional format on the difficulty table's
#        sku column also needs to follow the shift. ---
$ws.Range("C72:E74").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("C72:E74"))
